# Hortaliza, Terminal La Palmera de La Serena - Melón
# A new week of price observations (market date 2023-01-17, serial 44943)
# was inserted into the daily log, ahead of the existing 2022-12-28 (44923)
# block. This pushes all subsequent rows (old rows 92-118) down by three
# rows (new rows 95-121), growing the sheet from A1:R118 to A1:R121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 92, shifting the
# remainder of the table (old rows 92-118) down to 95-121.
$ws.Range("A92:R94").EntireRow.Insert()

$newRows = @(
    @{Row=92; Fecha=44943; Calidad="Extra";   Volumen=2400; Min=1800; Max=2000; Prom=1900; Origen="Región de O'Higgins"; PrecioKg=1900},
    @{Row=93; Fecha=44943; Calidad="Primera"; Volumen=2000; Min=1400; Max=1500; Prom=1450; Origen="Región de O'Higgins"; PrecioKg=1450},
    @{Row=94; Fecha=44943; Calidad="Segunda"; Volumen=1800; Min=1000; Max=1200; Prom=1100; Origen="Región de O'Higgins"; PrecioKg=1100}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 8
    $ws.Range("B$row").Value = "Terminal La Palmera de La Serena"
    $ws.Range("C$row").Value = "Coquimbo"
    $ws.Range("D$row").Value = $r.Fecha
    $ws.Range("E$row").Value = 4
    $ws.Range("F$row").Value = 100112027
    $ws.Range("G$row").Value = "Melón"
    $ws.Range("H$row").Value = "Tuna"
    $ws.Range("I$row").Value = $r.Calidad
    $ws.Range("J$row").Value = $r.Volumen
    $ws.Range("K$row").Value = $r.Min
    $ws.Range("L$row").Value = $r.Max
    $ws.Range("M$row").Value = $r.Prom
    $ws.Range("N$row").Value = "`$/unidad"
    $ws.Range("O$row").Value = $r.Origen
    $ws.Range("P$row").Value = $r.PrecioKg
    $ws.Range("Q$row").Value = 1
    $ws.Range("R$row").Value = "Hortaliza"
}
